$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 ("Rule" column, row 11 of the "R40" rule) is being renamed from
# "R40" to "1". In the source file this cell is stored as a shared string
# (text), so force a text-typed entry instead of letting "1" be
# auto-recognised as a number.
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "1"
